# Update cryptos list data per diff (Jan 15 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$origStyle_2 = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.652.90"
$ws.Range("D2").Style = $origStyle_2
$ws.Range("E2").Value = "  -0.37%  "

# Row 3
$origStyle_3 = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.542.67"
$ws.Range("D3").Style = $origStyle_3
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$origStyle_4 = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = $origStyle_4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$origStyle_5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.92"
$ws.Range("D5").Style = $origStyle_5
$ws.Range("E5").Value = "  +4.70%  "

# Row 6
$origStyle_6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.97"
$ws.Range("D6").Style = $origStyle_6
$ws.Range("E6").Value = "  -2.60%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$origStyle_9 = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.532"
$ws.Range("D9").Style = $origStyle_9
$ws.Range("E9").Value = "  -2.30%  "

# Row 10
$origStyle_10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.28"
$ws.Range("D10").Style = $origStyle_10
$ws.Range("E10").Value = "  -1.15%  "

# Row 11
$ws.Range("E11").Value = "  -1.16%  "

# Row 12
$origStyle_12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.59"
$ws.Range("D12").Style = $origStyle_12
$ws.Range("E12").Value = "  +0.51%  "

# Row 13
$ws.Range("E13").Value = "  -0.24%  "

# Row 14
$origStyle_14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.931.48"
$ws.Range("D14").Style = $origStyle_14
$ws.Range("E14").Value = "  +0.19%  "

# Row 15
$origStyle_15 = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.90"
$ws.Range("D15").Style = $origStyle_15
$ws.Range("E15").Value = "  +5.97%  "

# Row 16
$origStyle_16 = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.509.57"
$ws.Range("D16").Style = $origStyle_16
$ws.Range("E16").Value = "  -2.16%  "

# Row 17
$origStyle_17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.863"
$ws.Range("D17").Style = $origStyle_17
$ws.Range("E17").Value = "  -0.66%  "

# Row 18
$origStyle_18 = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.694.74"
$ws.Range("D18").Style = $origStyle_18
$ws.Range("E18").Value = "  -0.30%  "

# Row 19
$origStyle_19 = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.09"
$ws.Range("D19").Style = $origStyle_19
$ws.Range("E19").Value = "  -0.70%  "

# Row 20
$origStyle_20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.65"
$ws.Range("D20").Style = $origStyle_20
$ws.Range("E20").Value = "  +1.35%  "

# Row 21
$ws.Range("E21").Value = "  -1.77%  "

# Row 22
$origStyle_22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.01"
$ws.Range("D22").Style = $origStyle_22
$ws.Range("E22").Value = "  -0.76%  "

# Row 23
$origStyle_23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.04"
$ws.Range("D23").Style = $origStyle_23
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("E24").Value = "  +1.43%  "

# Row 25
$ws.Range("E25").Value = "  -2.58%  "

# Row 26
$origStyle_26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.20"
$ws.Range("D26").Style = $origStyle_26
$ws.Range("E26").Value = "  -2.09%  "

# Row 27
$origStyle_27 = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = $origStyle_27
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("E28").Value = "  +3.58%  "

# Row 29
$origStyle_29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.24"
$ws.Range("D29").Style = $origStyle_29
$ws.Range("E29").Value = "  +0.91%  "

# Row 30
$ws.Range("E30").Value = "  +3.89%  "

# Row 31
$ws.Range("E31").Value = "  -3.23%  "

# Row 32
$origStyle_32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.07"
$ws.Range("D32").Style = $origStyle_32
$ws.Range("E32").Value = "  -0.71%  "

# Row 33
$origStyle_33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("D33").Style = $origStyle_33
$ws.Range("E33").Value = "  +1.60%  "

# Row 34
$origStyle_34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.35"
$ws.Range("D34").Style = $origStyle_34
$ws.Range("E34").Value = "  +1.20%  "

# Row 35
$origStyle_35 = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.23"
$ws.Range("D35").Style = $origStyle_35
$ws.Range("E35").Value = "  -0.94%  "

# Row 36
$origStyle_36 = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0790"
$ws.Range("D36").Style = $origStyle_36
$ws.Range("E36").Value = "  -0.52%  "

# Row 37
$origStyle_37 = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.63"
$ws.Range("D37").Style = $origStyle_37
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("E38").Value = "  -2.83%  "

# Row 39
$ws.Range("B39").Value = "ApeXProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$origStyle_39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.44"
$ws.Range("D39").Style = $origStyle_39
$ws.Range("E39").Value = "  +13.32%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$origStyle_40 = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("D40").Style = $origStyle_40
$ws.Range("E40").Value = "  -0.40%  "

# Row 41
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$origStyle_41 = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.46"
$ws.Range("D41").Style = $origStyle_41
$ws.Range("E41").Value = "  -6.21%  "

# Row 42
$ws.Range("E42").Value = "  -0.40%  "

# Row 44
$ws.Range("E44").Value = "  +0.30%  "

# Row 45
$origStyle_45 = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0302"
$ws.Range("D45").Style = $origStyle_45
$ws.Range("E45").Value = "  -0.55%  "

# Row 46
$origStyle_46 = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.028.02"
$ws.Range("D46").Style = $origStyle_46
$ws.Range("E46").Value = "  -3.05%  "

# Row 47
$origStyle_47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "84.44"
$ws.Range("D47").Style = $origStyle_47
$ws.Range("E47").Value = "  -2.04%  "

# Row 48
$origStyle_48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.92"
$ws.Range("D48").Style = $origStyle_48
$ws.Range("E48").Value = "  +0.05%  "

# Row 49
$origStyle_49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.785.14"
$ws.Range("D49").Style = $origStyle_49
$ws.Range("E49").Value = "  -0.06%  "

# Row 50
$origStyle_50 = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.66"
$ws.Range("D50").Style = $origStyle_50
$ws.Range("E50").Value = "  +0.21%  "

# Row 51
$ws.Range("E51").Value = "  -0.50%  "
